$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 6) - "Roadsters" category, child of "Cars"
$ws.Range("A6").Value = "Roadsters"
$ws.Range("B6").Value = "Roadsters"
$ws.Range("C6").Value = "Roadsters"
$ws.Range("D6").Value = "Cars"

# Normalize the bold/duplicate style previously applied to column B and D
# (s="1") back to the default style (s="0") for rows 2-4.
$ws.Range("B2").NumberFormat = "General"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "General"

# Adjust column widths (B, C, D) to match the re-saved layout.
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(4).ColumnWidth = 8.0

# Move the active selection back to A1.
[void]$ws.Range("A1").Select()
